# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 3
    9  = 9
    10 = 3
    11 = 3
    12 = 4
    13 = 3
    14 = 4
    15 = 6
    16 = 1
    17 = 3
    18 = 3
    19 = 3
    20 = 3
    21 = 5
    22 = 3
    23 = 5
    24 = 5
    25 = 1
    27 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
